# Regenerate the localization-status report:
#  - flip the in-flight status from "Ready for handoff" to "In Translation"
#    (Overview!E2/F2 and the per-locale sheets' C2 all share this text)
#  - narrow the "zh-cn"/"de-de" status columns to match the refreshed layout

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E").ColumnWidth = 12.42
$overview.Columns("F").ColumnWidth = 12.42

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C").ColumnWidth = 12.42

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C").ColumnWidth = 12.42
